# Update the "Förändrad" (C) column date stamp for every existing data row
# (rows 2-170) from 2023-09-17 (45186) to 2023-09-19 (45188).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 170; $r++) {
    $ws.Cells.Item($r, 3).Value = 45188
}

# Row 170 picks up an explicit row height in the new file.
$ws.Rows.Item(170).RowHeight = 15

# Append the new record as row 171.
$ws.Cells.Item(171, 1).Value = "A 43784-2023"

$ws.Cells.Item(171, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(171, 2).Value = 45187

$ws.Cells.Item(171, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(171, 3).Value = 45188

$ws.Cells.Item(171, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(171, 5).Value = "GNOSJÖ"

$ws.Cells.Item(171, 7).Value = 3.8
$ws.Cells.Item(171, 8).Value = 0
$ws.Cells.Item(171, 9).Value = 0
$ws.Cells.Item(171, 10).Value = 0
$ws.Cells.Item(171, 11).Value = 0
$ws.Cells.Item(171, 12).Value = 0
$ws.Cells.Item(171, 13).Value = 0
$ws.Cells.Item(171, 14).Value = 0
$ws.Cells.Item(171, 15).Value = 0
$ws.Cells.Item(171, 16).Value = 0
$ws.Cells.Item(171, 17).Value = 0

# R column keeps the wrap-text style used throughout the sheet, with no text.
$ws.Cells.Item(171, 18).WrapText = $true
